$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original values for columns C and E, rows 4-9
$oldC4 = $ws.Cells.Item(4, 3).Value()
$oldC5 = $ws.Cells.Item(5, 3).Value()
$oldC6 = $ws.Cells.Item(6, 3).Value()
$oldC7 = $ws.Cells.Item(7, 3).Value()
$oldC8 = $ws.Cells.Item(8, 3).Value()
$oldC9 = $ws.Cells.Item(9, 3).Value()

$oldE4 = $ws.Cells.Item(4, 5).Value()
$oldE5 = $ws.Cells.Item(5, 5).Value()
$oldE6 = $ws.Cells.Item(6, 5).Value()
$oldE7 = $ws.Cells.Item(7, 5).Value()
$oldE8 = $ws.Cells.Item(8, 5).Value()
$oldE9 = $ws.Cells.Item(9, 5).Value()

# Apply a downward circular rotation: row r gets the value that was
# previously in row r-1 (row 4 wraps around and gets the old row 9 value)
$ws.Cells.Item(4, 3).Value = $oldC9
$ws.Cells.Item(5, 3).Value = $oldC4
$ws.Cells.Item(6, 3).Value = $oldC5
$ws.Cells.Item(7, 3).Value = $oldC6
$ws.Cells.Item(8, 3).Value = $oldC7
$ws.Cells.Item(9, 3).Value = $oldC8

$ws.Cells.Item(4, 5).Value = $oldE9
$ws.Cells.Item(5, 5).Value = $oldE4
$ws.Cells.Item(6, 5).Value = $oldE5
$ws.Cells.Item(7, 5).Value = $oldE6
$ws.Cells.Item(8, 5).Value = $oldE7
$ws.Cells.Item(9, 5).Value = $oldE8
